# Scheduled runner update: refresh computed profit/price figures across
# the per-job leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) with the
# latest market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1887.5
$ws.Range("I100").Value = 1820
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1820
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1279
$ws.Range("N100").Value = -3082

$ws.Range("H107").Value = 951.85
$ws.Range("I107").Value = 951.85
$ws.Range("K107").Value = 951.85
$ws.Range("M107").Value = 968.15

$ws.Range("H112").Value = 76924450
$ws.Range("J112").Value = 76924450
$ws.Range("L112").Value = 230773350
$ws.Range("N112").Value = -230775566

$ws.Range("H113").Value = 2160.8
$ws.Range("J113").Value = 2037.4546
$ws.Range("L113").Value = 2037.4546
$ws.Range("N113").Value = -8545.454600000001

$ws.Range("H129").Value = 1246.8
$ws.Range("J129").Value = 1443.3334
$ws.Range("L129").Value = 4330.0002
$ws.Range("N129").Value = -14330.0002

$ws.Range("H132").Value = 326101.2
$ws.Range("I132").Value = 374257.94
$ws.Range("J132").Value = 1043.25
$ws.Range("K132").Value = 1122773.82
$ws.Range("L132").Value = 3129.75
$ws.Range("M132").Value = -1120243.82
$ws.Range("N132").Value = -8189.75

$ws.Range("H137").Value = 4397.275
$ws.Range("I137").Value = 5932.4546
$ws.Range("K137").Value = 17797.3638
$ws.Range("M137").Value = -15247.3638

$ws.Range("H138").Value = 1432397.1
$ws.Range("I138").Value = 2186.9583
$ws.Range("J138").Value = 2752591
$ws.Range("K138").Value = 6560.874899999999
$ws.Range("L138").Value = 8257773
$ws.Range("M138").Value = -1420.874899999999
$ws.Range("N138").Value = -8268053

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5882802.5
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = 7353366
$ws.Range("K2").Value = 550
$ws.Range("L2").Value = 7353366
$ws.Range("M2").Value = -437
$ws.Range("N2").Value = -7353592

$ws.Range("H32").Value = 7803.653
$ws.Range("I32").Value = 5022.698
$ws.Range("J32").Value = 15561.053
$ws.Range("K32").Value = 5022.698
$ws.Range("L32").Value = 15561.053
$ws.Range("M32").Value = -4735.698
$ws.Range("N32").Value = -16135.053

$ws.Range("H45").Value = 49922.047
$ws.Range("I45").Value = 85104.086
$ws.Range("J45").Value = 3012.6667
$ws.Range("K45").Value = 85104.086
$ws.Range("L45").Value = 3012.6667
$ws.Range("M45").Value = -84727.086
$ws.Range("N45").Value = -3766.6667

$ws.Range("H61").Value = 2706.0715
$ws.Range("I61").Value = 2398.9048
$ws.Range("J61").Value = 3627.5715
$ws.Range("K61").Value = 2398.9048
$ws.Range("L61").Value = 3627.5715
$ws.Range("M61").Value = -2186.9048
$ws.Range("N61").Value = -4051.5715

$ws.Range("H109").Value = 34458.668
$ws.Range("J109").Value = 34458.668
$ws.Range("L109").Value = 34458.668
$ws.Range("N109").Value = -37232.668

$ws.Range("H110").Value = 58031.07
$ws.Range("I110").Value = 921.36365
$ws.Range("J110").Value = 267433.34
$ws.Range("K110").Value = 921.36365
$ws.Range("L110").Value = 267433.34
$ws.Range("M110").Value = 1123.63635
$ws.Range("N110").Value = -271523.34

$ws.Range("H116").Value = 5882802.5
$ws.Range("I116").Value = 550
$ws.Range("J116").Value = 7353366
$ws.Range("K116").Value = 550
$ws.Range("L116").Value = 7353366
$ws.Range("M116").Value = 1744
$ws.Range("N116").Value = -7357954

$ws.Range("H133").Value = 99260.25
$ws.Range("J133").Value = 99260.25
$ws.Range("L133").Value = 99260.25
$ws.Range("N133").Value = -104320.25

$ws.Range("H136").Value = 2706.0715
$ws.Range("I136").Value = 2398.9048
$ws.Range("J136").Value = 3627.5715
$ws.Range("K136").Value = 7196.714399999999
$ws.Range("L136").Value = 10882.7145
$ws.Range("M136").Value = -4646.714399999999
$ws.Range("N136").Value = -15982.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5882802.5
$ws.Range("I3").Value = 550
$ws.Range("J3").Value = 7353366
$ws.Range("K3").Value = 550
$ws.Range("L3").Value = 7353366
$ws.Range("M3").Value = -436
$ws.Range("N3").Value = -7353594

$ws.Range("H99").Value = 1177.7778
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 775
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 775
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -3771

$ws.Range("H107").Value = 614674.7
$ws.Range("I107").Value = 881888.75
$ws.Range("J107").Value = 3899.7144
$ws.Range("K107").Value = 881888.75
$ws.Range("L107").Value = 3899.7144
$ws.Range("M107").Value = -879968.75
$ws.Range("N107").Value = -7739.7144

$ws.Range("H134").Value = 45816.64
$ws.Range("I134").Value = 59185.05
$ws.Range("K134").Value = 177555.15
$ws.Range("M134").Value = -175020.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2303
$ws.Range("I16").Value = 1718.5714
$ws.Range("J16").Value = 3666.6667
$ws.Range("K16").Value = 1718.5714
$ws.Range("L16").Value = 3666.6667
$ws.Range("M16").Value = -1431.5714
$ws.Range("N16").Value = -4240.6667

$ws.Range("H31").Value = 1743.275
$ws.Range("I31").Value = 1228.5
$ws.Range("K31").Value = 1228.5
$ws.Range("M31").Value = -933.5

$ws.Range("H34").Value = 1743.275
$ws.Range("I34").Value = 1228.5
$ws.Range("K34").Value = 1228.5
$ws.Range("M34").Value = -1026.5

$ws.Range("H107").Value = 3501.8667
$ws.Range("I107").Value = 3612.75
$ws.Range("J107").Value = 3375.1428
$ws.Range("K107").Value = 3612.75
$ws.Range("L107").Value = 3375.1428
$ws.Range("M107").Value = -1692.75
$ws.Range("N107").Value = -7215.1428

$ws.Range("H113").Value = 2303
$ws.Range("I113").Value = 1718.5714
$ws.Range("J113").Value = 3666.6667
$ws.Range("K113").Value = 1718.5714
$ws.Range("L113").Value = 3666.6667
$ws.Range("M113").Value = 451.4286
$ws.Range("N113").Value = -8006.6667

$ws.Range("H134").Value = 2325.484
$ws.Range("I134").Value = 2372.739
$ws.Range("J134").Value = 2189.625
$ws.Range("K134").Value = 7118.217000000001
$ws.Range("L134").Value = 6568.875
$ws.Range("M134").Value = -4583.217000000001
$ws.Range("N134").Value = -11638.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 836.63855
$ws.Range("I68").Value = 610.12964
$ws.Range("J68").Value = 1258.4138
$ws.Range("K68").Value = 1830.38892
$ws.Range("L68").Value = 3775.2414
$ws.Range("M68").Value = -1019.38892
$ws.Range("N68").Value = -5397.2414

$ws.Range("H71").Value = 836.63855
$ws.Range("I71").Value = 610.12964
$ws.Range("J71").Value = 1258.4138
$ws.Range("K71").Value = 5491.16676
$ws.Range("L71").Value = 11325.7242
$ws.Range("M71").Value = -1435.16676
$ws.Range("N71").Value = -19437.7242

$ws.Range("H87").Value = 12000
$ws.Range("I87").Value = 5000
$ws.Range("J87").Value = 13000
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 39000
$ws.Range("M87").Value = -13752
$ws.Range("N87").Value = -41496

$ws.Range("H90").Value = 12000
$ws.Range("I90").Value = 5000
$ws.Range("J90").Value = 13000
$ws.Range("K90").Value = 45000
$ws.Range("L90").Value = 117000
$ws.Range("M90").Value = -38760
$ws.Range("N90").Value = -129480

$ws.Range("H107").Value = 38561.133
$ws.Range("I107").Value = 27643.838
$ws.Range("J107").Value = 63807.375
$ws.Range("K107").Value = 82931.514
$ws.Range("L107").Value = 191422.125
$ws.Range("M107").Value = -81011.514
$ws.Range("N107").Value = -195262.125

$ws.Range("H131").Value = 1962937
$ws.Range("J131").Value = 2705278
$ws.Range("L131").Value = 8115834
$ws.Range("N131").Value = -8125914

$ws.Range("H132").Value = 71429830
$ws.Range("I132").Value = 142858780
$ws.Range("J132").Value = 885.1429000000001
$ws.Range("K132").Value = 1285729020
$ws.Range("L132").Value = 7966.2861
$ws.Range("M132").Value = -1285726490
$ws.Range("N132").Value = -13026.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1332.36
$ws.Range("I102").Value = 1295.7142
$ws.Range("J102").Value = 1379
$ws.Range("K102").Value = 1295.7142
$ws.Range("L102").Value = 1379
$ws.Range("M102").Value = 326.2858000000001
$ws.Range("N102").Value = -4623

$ws.Range("H122").Value = 4792.864
$ws.Range("J122").Value = 3018
$ws.Range("L122").Value = 9054
$ws.Range("N122").Value = -13954

$ws.Range("H126").Value = 1687.5
$ws.Range("I126").Value = 1375
$ws.Range("K126").Value = 4125
$ws.Range("M126").Value = -1655

$ws.Range("H132").Value = 2562
$ws.Range("I132").Value = 3832.2
$ws.Range("J132").Value = 1715.2
$ws.Range("K132").Value = 11496.6
$ws.Range("L132").Value = 5145.6
$ws.Range("M132").Value = -8966.599999999999
$ws.Range("N132").Value = -10205.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5509.3335
$ws.Range("I40").Value = 5940.5713
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 5940.5713
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -5804.5713
$ws.Range("N40").Value = -4272

$ws.Range("H133").Value = 47000
$ws.Range("J133").Value = 47000
$ws.Range("L133").Value = 47000
$ws.Range("N133").Value = -52060

$ws.Range("H134").Value = 46412.23
$ws.Range("J134").Value = 46412.23
$ws.Range("L134").Value = 46412.23
$ws.Range("N134").Value = -56552.23

